# Add 2022-Q3 data
# 1) Update the "总计" (totals) summary sheet: insert a new row for 2022-Q3
#    right under the header, shifting all existing rows down by one.
# 2) Insert a brand-new "2022-Q3" worksheet (positioned right after "总计",
#    before "2022-Q2") containing the per-fund holdings data for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" sheet - insert new row 2 for 2022-Q3
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

$totals.Rows.Item(2).Insert()
$totals.Range("A2:D2").ClearFormats()

# Copy column-A cell style (bold + border + centered) from the row below,
# which already carries the correct formatting (same sheet copy keeps style).
$totals.Cells.Item(3,1).Copy($totals.Cells.Item(2,1))

$totals.Cells.Item(2,1).Value = 0
$totals.Cells.Item(2,2).Value = "2022-Q3"
$totals.Cells.Item(2,3).Value = 3
$totals.Cells.Item(2,4).Value = 0.11

# The A column holds a 0-based running index; renumber the rows that were
# shifted down so the sequence stays 0,1,2,3,4,5,6.
$totals.Cells.Item(3,1).Value = 1
$totals.Cells.Item(4,1).Value = 2
$totals.Cells.Item(5,1).Value = 3
$totals.Cells.Item(6,1).Value = 4
$totals.Cells.Item(7,1).Value = 5
$totals.Cells.Item(8,1).Value = 6

# ---------------------------------------------------------------------
# Step 2: create the "2022-Q3" worksheet by duplicating "2022-Q2" (so the
# formatting / styles / text-vs-number cell types are preserved exactly),
# then overwrite the data with the 2022-Q3 figures, and trim extra rows.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Remove the now-duplicate rows 5-7 (2022-Q3 only has 3 fund rows, vs the 5
# that 2022-Q2 had), shifting remaining cells up.
$q3.Range("A5:H7").Delete(-4162)

# Columns B,C,D,E,F,G must stay text (fund code / decimal-formatted strings
# stored as text, matching original source data), so force text format
# before writing, then clear the temporary number format back off again.
$textRange = $q3.Range("B2:G4")
$textRange.NumberFormat = "@"

$q3.Cells.Item(2,2).Value = "000849"
$q3.Cells.Item(2,3).Value = "汇丰晋信双核策略混合A"
$q3.Cells.Item(2,4).Value = "2.11"
$q3.Cells.Item(2,5).Value = "64.43"
$q3.Cells.Item(2,6).Value = "3.78"
$q3.Cells.Item(2,7).Value = "0.0798"
$q3.Cells.Item(2,8).Value = 7

$q3.Cells.Item(3,2).Value = "510200"
$q3.Cells.Item(3,3).Value = "汇安上证证券ETF"
$q3.Cells.Item(3,4).Value = "0.67"
$q3.Cells.Item(3,5).Value = "95.06"
$q3.Cells.Item(3,6).Value = "3.06"
$q3.Cells.Item(3,7).Value = "0.0205"
$q3.Cells.Item(3,8).Value = 10

$q3.Cells.Item(4,2).Value = "000850"
$q3.Cells.Item(4,3).Value = "汇丰晋信双核策略混合C"
$q3.Cells.Item(4,4).Value = "0.33"
$q3.Cells.Item(4,5).Value = "64.43"
$q3.Cells.Item(4,6).Value = "3.78"
$q3.Cells.Item(4,7).Value = "0.0125"
$q3.Cells.Item(4,8).Value = 7

$textRange.ClearFormats()
